$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2718.0667
$ws.Range("I9").Value = 406.85715
$ws.Range("K9").Value = 406.85715
$ws.Range("M9").Value = -237.85715
$ws.Range("H19").Value = 226
$ws.Range("I19").Value = 166.4
$ws.Range("K19").Value = 166.4
$ws.Range("M19").Value = 8.599999999999994
$ws.Range("H32").Value = 1608.4286
$ws.Range("I32").Value = 1649.8
$ws.Range("K32").Value = 1649.8
$ws.Range("M32").Value = -1323.8
$ws.Range("H51").Value = 12500
$ws.Range("I51").Value = 12500
$ws.Range("K51").Value = 12500
$ws.Range("M51").Value = -12016
$ws.Range("H69").Value = 5500
$ws.Range("H72").Value = 5500
$ws.Range("H74").Value = 3783.3333
$ws.Range("I74").Value = 3783.3333
$ws.Range("K74").Value = 3783.3333
$ws.Range("M74").Value = -2847.3333
$ws.Range("H77").Value = 3783.3333
$ws.Range("I77").Value = 3783.3333
$ws.Range("K77").Value = 18916.6665
$ws.Range("M77").Value = -14236.6665
$ws.Range("H112").Value = 2440.5
$ws.Range("I112").Value = 1323.25
$ws.Range("J112").Value = 2999.125
$ws.Range("K112").Value = 3969.75
$ws.Range("L112").Value = 8997.375
$ws.Range("M112").Value = -2861.75
$ws.Range("N112").Value = -11213.375
$ws.Range("H116").Value = 6007.676
$ws.Range("I116").Value = 5281.4707
$ws.Range("K116").Value = 5281.4707
$ws.Range("M116").Value = -1839.4707
$ws.Range("H137").Value = 3562.5
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 3928.5715
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 11785.7145
$ws.Range("M137").Value = -450
$ws.Range("N137").Value = -16885.7145
$ws.Range("H138").Value = 8193.5625
$ws.Range("I138").Value = 6273.1816
$ws.Range("K138").Value = 18819.5448
$ws.Range("M138").Value = -13679.5448

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5129.1914
$ws.Range("I32").Value = 4023.8445
$ws.Range("K32").Value = 4023.8445
$ws.Range("M32").Value = -3736.8445
$ws.Range("H61").Value = 3421
$ws.Range("J61").Value = 3426.25
$ws.Range("L61").Value = 3426.25
$ws.Range("N61").Value = -3850.25
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240
$ws.Range("H82").Value = 60090.5
$ws.Range("J82").Value = 80181
$ws.Range("L82").Value = 80181
$ws.Range("N82").Value = -80903
$ws.Range("H85").Value = 60090.5
$ws.Range("J85").Value = 80181
$ws.Range("L85").Value = 80181
$ws.Range("N85").Value = -82677
$ws.Range("H132").Value = 1473.3
$ws.Range("I132").Value = 637
$ws.Range("K132").Value = 1911
$ws.Range("M132").Value = 619
$ws.Range("H136").Value = 3421
$ws.Range("J136").Value = 3426.25
$ws.Range("L136").Value = 10278.75
$ws.Range("N136").Value = -15378.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14291.4
$ws.Range("I82").Value = 14291.4
$ws.Range("K82").Value = 14291.4
$ws.Range("M82").Value = -13908.4
$ws.Range("H85").Value = 14291.4
$ws.Range("I85").Value = 14291.4
$ws.Range("K85").Value = 14291.4
$ws.Range("M85").Value = -12965.4
$ws.Range("H86").Value = 3470.125
$ws.Range("I86").Value = 3118
$ws.Range("J86").Value = 3681.4
$ws.Range("K86").Value = 3118
$ws.Range("L86").Value = 3681.4
$ws.Range("M86").Value = -1995
$ws.Range("N86").Value = -5927.4
$ws.Range("H89").Value = 3470.125
$ws.Range("I89").Value = 3118
$ws.Range("J89").Value = 3681.4
$ws.Range("K89").Value = 15590
$ws.Range("L89").Value = 18407
$ws.Range("M89").Value = -9974
$ws.Range("N89").Value = -29639
$ws.Range("H134").Value = 1852.3158
$ws.Range("I134").Value = 1683.5625
$ws.Range("K134").Value = 5050.6875
$ws.Range("M134").Value = -2515.6875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H132").Value = 6792
$ws.Range("I132").Value = 2688
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 8064
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -5534
$ws.Range("N132").Value = -50060
$ws.Range("H134").Value = 3793.4167
$ws.Range("J134").Value = 4999
$ws.Range("L134").Value = 14997
$ws.Range("N134").Value = -20067
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 270.16666
$ws.Range("I11").Value = 340.33334
$ws.Range("K11").Value = 1021.00002
$ws.Range("M11").Value = -881.0000200000001
$ws.Range("H24").Value = 816.3333
$ws.Range("I24").Value = 621.4286
$ws.Range("J24").Value = 1498.5
$ws.Range("K24").Value = 1864.2858
$ws.Range("L24").Value = 4495.5
$ws.Range("M24").Value = -1634.2858
$ws.Range("N24").Value = -4955.5
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H92").Value = 445.125
$ws.Range("I92").Value = 282.2
$ws.Range("J92").Value = 716.6667
$ws.Range("K92").Value = 846.5999999999999
$ws.Range("L92").Value = 2150.0001
$ws.Range("M92").Value = 401.4000000000001
$ws.Range("N92").Value = -4646.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1084532.6
$ws.Range("I11").Value = 1084532.6
$ws.Range("K11").Value = 1084532.6
$ws.Range("M11").Value = -1084393.6
$ws.Range("H43").Value = 2606
$ws.Range("J43").Value = 2019
$ws.Range("L43").Value = 2019
$ws.Range("N43").Value = -2321
$ws.Range("H117").Value = 67897.5
$ws.Range("J117").Value = 67897.5
$ws.Range("L117").Value = 67897.5
$ws.Range("N117").Value = -74781.5
$ws.Range("H122").Value = 145996.86
$ws.Range("I122").Value = 202799.2
$ws.Range("J122").Value = 3991
$ws.Range("K122").Value = 608397.6000000001
$ws.Range("L122").Value = 11973
$ws.Range("M122").Value = -605947.6000000001
$ws.Range("N122").Value = -16873
$ws.Range("H126").Value = 200002000
$ws.Range("I126").Value = 250002000
$ws.Range("K126").Value = 750006000
$ws.Range("M126").Value = -750003530
$ws.Range("H132").Value = 3468.889
$ws.Range("I132").Value = 2870.1667
$ws.Range("K132").Value = 8610.500100000001
$ws.Range("M132").Value = -6080.500100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 16012.667
$ws.Range("J17").Value = 22069
$ws.Range("L17").Value = 22069
$ws.Range("N17").Value = -22409
$ws.Range("H39").Value = 8000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 2258.4
$ws.Range("H132").Value = 7470.857
$ws.Range("I132").Value = 7049.3335
$ws.Range("K132").Value = 21148.0005
$ws.Range("M132").Value = -18618.0005
$ws.Range("H136").Value = 4498.5
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
